# Prepared Spreadsheet for Devlery
# Updated the Topical Areas to match current 3.10 categories.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Find-RowByLabel([string]$label) {
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
    for ($r = 1; $r -le $lastRow; $r++) {
        $v = $ws.Cells.Item($r, 1).Value2
        if ($v -eq $label) {
            return $r
        }
    }
    return -1
}

# 1. Remove the "Business" topical area row entirely.
$rowBusiness = Find-RowByLabel "Business"
if ($rowBusiness -gt 0) {
    $ws.Rows.Item($rowBusiness).Delete()
}

# 2. Remove the "Teen Crafts" topical area row entirely.
$rowTeenCrafts = Find-RowByLabel "Teen Crafts"
if ($rowTeenCrafts -gt 0) {
    $ws.Rows.Item($rowTeenCrafts).Delete()
}

# 3. Insert the new "series/African American History Month" row directly
#    above the existing "series/Artober" row.
$rowArtober = Find-RowByLabel "series/Artober"
if ($rowArtober -gt 0) {
    $ws.Rows.Item($rowArtober).Insert()
    $ws.Cells.Item($rowArtober, 1).Value = "series/African American History Month"
    $ws.Cells.Item($rowArtober, 3).Value = "X-BEDEWORK-ALIAS;X-BEDEWORK-PARAM-DISPLAYNAME=series/African American History Month:/user/agrp_calsuite-MainCampus/series/African American History Month"
}

# 4. Insert the new "series/Studio NPL" row directly above the existing
#    "series/Seed Exchange" row.
$rowSeedExchange = Find-RowByLabel "series/Seed Exchange"
if ($rowSeedExchange -gt 0) {
    $ws.Rows.Item($rowSeedExchange).Insert()
    $ws.Cells.Item($rowSeedExchange, 1).Value = "series/Studio NPL"
    $ws.Cells.Item($rowSeedExchange, 3).Value = "X-BEDEWORK-ALIAS;X-BEDEWORK-PARAM-DISPLAYNAME=series/Studio NPL:/user/agrp_calsuite-MainCampus/series/Studio NPL"
}

# Reset the view: select the full data range and scroll back to the top,
# replacing the old scrolled-down / single-cell selection state.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$ws.Range("A1:C" + $lastRow).Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
